$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Friday (row 17) worked 10:00 AM - 2:00 PM (4 hours)
$ws.Range("C17").Value = 0.416666666666667
$ws.Range("D17").Value = 0.583333333333333

# Remove the now-stale "10-11, " note that was next to the Friday row
$ws.Range("L17").Value = ""
$ws.Range("L17").NumberFormat = "General"

# Update the active selection to reflect where the user left off
$ws.Range("D18").Select()
